# Regenerate merged AHB files
#
# The "Änderung" (change) column L previously flagged a fixed set of rows
# with the literal text "ÄNDERUNG" (shared string 191, bold/gold style 7,
# centered). This edit clears that flag for a batch of rows: column L goes
# back to its normal empty/centered style (style 4, same as the untouched
# rows already use), and for a subset of "group header" rows the rest of
# the row (A:V) is also normalized from the bordered-only style (5/"B" bold)
# to the grey-filled style used by the other header rows (2, with B as 3).
#
# We reproduce this purely through formatting operations (copy/paste of
# formats from a template row/cell that already carries the desired style)
# so the engine reuses the existing cellXfs entries instead of minting new
# ones, then clear the stale "ÄNDERUNG" text out of column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose entire A:V span is restyled (group header rows) -- these also
# lose their column-L "ÄNDERUNG" flag.
$headerRows = @(30, 33, 40, 47, 51, 54, 78)

# Rows where only column L's "ÄNDERUNG" flag/style is cleared.
$lOnlyRows = @(31, 32, 34, 35, 36, 37, 38, 39, 41, 42, 44, 45, 48, 49, 50, 52, 53, 55, 56, 64, 65, 66, 72, 73, 74, 76, 77, 79, 80)

# Template row 2 already carries the exact target styles for every column
# (A=2, B=3, C:K=2, L=4, M:V=2), so copying its formats across reproduces
# the target cellXfs indices without creating new style entries.
$templateRow = $ws.Range("A2:V2")

foreach ($r in $headerRows) {
    $templateRow.Copy()
    $targetRange = "A" + $r + ":V" + $r
    $ws.Range($targetRange).PasteSpecial(-4122)
    $ws.Range("L" + $r).ClearContents()
}

$templateL = $ws.Range("L2")

foreach ($r in $lOnlyRows) {
    $templateL.Copy()
    $ws.Range("L" + $r).PasteSpecial(-4122)
    $ws.Range("L" + $r).ClearContents()
}

$excel.CutCopyMode = 0
